# Update countries & provincias Spain
# This script applies a refreshed data scrape to the "Pais" worksheet:
#  - Updates the "last updated" timestamp
#  - Refreshes several countries' case statistics
#  - Re-sorts a couple of countries (Guatemala, Butan) ahead of their
#    neighbours, shifting the rows in between down by one position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 15 de Mayo de 2020 a las 04:05"

# Brasil (row 9)
$ws.Range("B9").Value = 203165
$ws.Range("C9").Value = 247
$ws.Range("E9").Value = 109687
$ws.Range("G9").Value = 6
$ws.Range("H9").Value = 13999

# Corea del Sur (row 46)
$ws.Range("B46").Value = 11018
$ws.Range("C46").Value = 27
$ws.Range("D46").Value = 9821
$ws.Range("E46").Value = 937

# Panama (row 50)
$ws.Range("F50").Value = 72

# Guatemala moves up, now row 88, with refreshed totals
$ws.Range("A88").Value = "Guatemala"
$ws.Range("B88").Value = 1518
$ws.Range("C88").Value = 176
$ws.Range("D88").Value = 129
$ws.Range("E88").Value = 1360
$ws.Range("F88").Value = 5
$ws.Range("H88").Value = 29

# Lituania, Nueva Zelanda, Eslovaquia, Eslovenia each shift down one row
$ws.Range("A89").Value = "Lituania"
$ws.Range("B89").Value = 1511
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 934
$ws.Range("E89").Value = 523
$ws.Range("F89").Value = 17
$ws.Range("H89").Value = 54

$ws.Range("A90").Value = "Nueva Zelanda"
$ws.Range("B90").Value = 1498
$ws.Range("C90").Value = 1
$ws.Range("D90").Value = 1421
$ws.Range("E90").Value = 56
$ws.Range("F90").Value = 2
$ws.Range("H90").Value = 21

$ws.Range("A91").Value = "Eslovaquia"
$ws.Range("B91").Value = 1477
$ws.Range("D91").Value = 1112
$ws.Range("E91").Value = 338
$ws.Range("F91").Value = 5
$ws.Range("H91").Value = 27

$ws.Range("A92").Value = "Eslovenia"
$ws.Range("B92").Value = 1464
$ws.Range("D92").Value = 267
$ws.Range("E92").Value = 1094
$ws.Range("F92").Value = 7
$ws.Range("H92").Value = 103

# Butan moves up ahead of Mauritania, with refreshed totals
$ws.Range("A190").Value = "Butan"
$ws.Range("D190").Value = 5
$ws.Range("E190").Value = 15
$ws.Range("H190").Value = 0

# Mauritania shifts down one row
$ws.Range("A191").Value = "Mauritania"
$ws.Range("B191").Value = 20
$ws.Range("D191").Value = 7
$ws.Range("E191").Value = 11
$ws.Range("H191").Value = 2
